$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

# Force text format first on D-column cells whose new value would otherwise be
# auto-coerced to a Number (losing trailing zeros / changing cell type).
foreach ($addr in @("D5", "D6", "D10", "D14", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D32", "D35", "D36", "D37", "D38", "D39", "D40", "D46", "D47", "D48", "D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.299.03"
$ws.Range("E2").Value = "  +1.49%  "

$ws.Range("D3").Value = "3.137.36"

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "580.32"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").Value = "174.93"
$ws.Range("E6").Value = "  +3.69%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.133.04"
$ws.Range("E8").Value = "  +3.40%  "

$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("D10").Value = "6.50"
$ws.Range("E10").Value = "  -2.50%  "

$ws.Range("E11").Value = "  +2.15%  "

$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("D14").Value = "37.48"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("D16").Value = "3.654.92"
$ws.Range("E16").Value = "  +3.42%  "

$ws.Range("D17").Value = "67.305.31"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").Value = "3.136.56"
$ws.Range("E19").Value = "  +3.56%  "

$ws.Range("D20").Value = "16.15"
$ws.Range("E20").Value = "  -2.10%  "

$ws.Range("D21").Value = "488.26"
$ws.Range("E21").Value = "  +4.35%  "

$ws.Range("E22").Value = "  +0.93%  "

$ws.Range("D23").Value = "7.70"
$ws.Range("E23").Value = "  +4.05%  "

$ws.Range("D24").Value = "84.30"
$ws.Range("E24").Value = "  +1.46%  "

$ws.Range("D25").Value = "13.26"
$ws.Range("E25").Value = "  +3.92%  "

$ws.Range("D26").Value = "2.34"
$ws.Range("E26").Value = "  +3.27%  "

$ws.Range("D27").Value = "10.09"
$ws.Range("E27").Value = "  +0.50%  "

$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("E29").Value = "  -2.86%  "

$ws.Range("D30").Value = "2.40"
$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("D31").Value = "2.70"
$ws.Range("E31").Value = "  +1.61%  "

$ws.Range("D32").Value = "28.88"
$ws.Range("E32").Value = "  +2.40%  "

$ws.Range("D33").Value = "0.0${sub3}0997"
$ws.Range("E33").Value = "  +0.10%  "

$ws.Range("E34").Value = "  -3.47%  "

$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").Value = "5.94"
$ws.Range("E36").Value = "  +1.37%  "

$ws.Range("D37").Value = "0.990"
$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").Value = "47.55"
$ws.Range("E38").Value = "  -1.32%  "

$ws.Range("D39").Value = "2.12"
$ws.Range("E39").Value = "  +2.70%  "

$ws.Range("D40").Value = "50.13"
$ws.Range("E40").Value = "  +1.24%  "

$ws.Range("E41").Value = "  +0.27%  "

$ws.Range("E42").Value = "  +1.58%  "

$ws.Range("E43").Value = "  +0.17%  "

$ws.Range("E44").Value = "  -1.51%  "

$ws.Range("D45").Value = "2.845.63"
$ws.Range("E45").Value = "  +5.16%  "

$ws.Range("D48").Value = "135.93"
$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").Value = "24.96"
$ws.Range("E50").Value = "  +1.81%  "

$ws.Range("E51").Value = "  -0.35%  "

# Rows 46/47: VeChain and Bittensor swap ranking positions
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "0.0359"
$ws.Range("E46").Value = "  -0.53%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "384.39"
$ws.Range("E47").Value = "  +1.33%  "
